$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Child")

$values = @{
    2  = "7,8"
    3  = "-8,5"
    4  = "-8,-8"
    5  = "6,-2"
    6  = "-7,-7"
    7  = "0,-10"
    8  = "-4,-8"
    9  = "-5,-9"
    10 = "0,8"
    11 = "-1,-1"
    12 = "-2,-8"
    13 = "6,5"
    14 = "-6,0"
    15 = "3,9"
    16 = "-4,2"
    17 = "7,-8"
    18 = "-4,3"
    19 = "8,-5"
    20 = "-4,3"
    21 = "-1,-5"
}

foreach ($row in $values.Keys) {
    $ws.Range("D$row").Value = $values[$row]
}
